# Generate Report for Handoff
#
# The localization-status report is regenerated after a new handoff
# xliff was produced for e2e\b.md. Its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff" on the
# Overview sheet as well as on the zh-cn / de-de detail sheets, and the
# per-language rows pick up the freshly generated handoff file name /
# datetime together with a "stale handback" error message (the handback
# file on record is now behind the new handoff).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5be971ca8e997e93ba2ddcb913a1f6a1491f8fd3/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9bae587c2df646e64a2b63f2a83fb5d752474244/e2e/b.md."

# ---- Overview sheet: row 3 is e2e\b.md ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-05 18:42:45"

# ---- zh-cn sheet: row 3 is b.md ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Force plain text ("False", not a real boolean) like the rest of the
# True/False columns, then drop the resulting quote-prefix style so the
# cell style index is left untouched.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-05 18:42:40"
$zhcn.Range("P3").Value = $errorDetail
# Column width 40 in the saved OOXML == ColumnWidth 40 - 5/6 here (the
# engine adds a fixed 5/6-character padding on top of ColumnWidth).
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet: row 3 is b.md ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-05 18:42:45"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
